# Auto-generated edit script: updates market price / profit values
# across multiple sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(13, 8).Value = 5854.75  # H13: 4998.25 -> 5854.75
$ws.Cells.Item(13, 9).Value = 1000  # I13: 995 -> 1000
$ws.Cells.Item(13, 10).Value = 7473  # J13: 6332.6665 -> 7473
$ws.Cells.Item(13, 11).Value = 1000  # K13: 995 -> 1000
$ws.Cells.Item(13, 12).Value = 7473  # L13: 6332.6665 -> 7473
$ws.Cells.Item(13, 13).Value = -831  # M13: -826 -> -831
$ws.Cells.Item(13, 14).Value = -7811  # N13: -6670.6665 -> -7811
$ws.Cells.Item(18, 8).Value = 300  # H18: 891.5 -> 300
$ws.Cells.Item(18, 9).Value = 300  # I18: 891.5 -> 300
$ws.Cells.Item(18, 11).Value = 300  # K18: 891.5 -> 300
$ws.Cells.Item(18, 13).Value = -16  # M18: -607.5 -> -16
$ws.Cells.Item(40, 8).Value = 1160.2903  # H40: 1144.1515 -> 1160.2903
$ws.Cells.Item(40, 10).Value = 1249.9166  # J40: 1199.0714 -> 1249.9166
$ws.Cells.Item(40, 12).Value = 1249.9166  # L40: 1199.0714 -> 1249.9166
$ws.Cells.Item(40, 14).Value = -1599.9166  # N40: -1549.0714 -> -1599.9166
$ws.Cells.Item(64, 8).Value = 4817  # H64: 5283.5 -> 4817
$ws.Cells.Item(64, 9).Value = 5563.3335  # I64: 6720 -> 5563.3335
$ws.Cells.Item(64, 10).Value = 3697.5  # J64: 3847 -> 3697.5
$ws.Cells.Item(64, 11).Value = 5563.3335  # K64: 6720 -> 5563.3335
$ws.Cells.Item(64, 12).Value = 3697.5  # L64: 3847 -> 3697.5
$ws.Cells.Item(64, 13).Value = -5315.3335  # M64: -6472 -> -5315.3335
$ws.Cells.Item(64, 14).Value = -4193.5  # N64: -4343 -> -4193.5
$ws.Cells.Item(67, 8).Value = 4817  # H67: 5283.5 -> 4817
$ws.Cells.Item(67, 9).Value = 5563.3335  # I67: 6720 -> 5563.3335
$ws.Cells.Item(67, 10).Value = 3697.5  # J67: 3847 -> 3697.5
$ws.Cells.Item(67, 11).Value = 5563.3335  # K67: 6720 -> 5563.3335
$ws.Cells.Item(67, 12).Value = 3697.5  # L67: 3847 -> 3697.5
$ws.Cells.Item(67, 13).Value = -4705.3335  # M67: -5862 -> -4705.3335
$ws.Cells.Item(67, 14).Value = -5413.5  # N67: -5563 -> -5413.5
$ws.Cells.Item(86, 8).Value = 3004.5  # H86: 3019.45 -> 3004.5
$ws.Cells.Item(86, 9).Value = 2446.0667  # I86: 2368.125 -> 2446.0667
$ws.Cells.Item(86, 10).Value = 4679.8  # J86: 5624.75 -> 4679.8
$ws.Cells.Item(86, 11).Value = 2446.0667  # K86: 2368.125 -> 2446.0667
$ws.Cells.Item(86, 12).Value = 4679.8  # L86: 5624.75 -> 4679.8
$ws.Cells.Item(86, 13).Value = -1323.0667  # M86: -1245.125 -> -1323.0667
$ws.Cells.Item(86, 14).Value = -6925.8  # N86: -7870.75 -> -6925.8
$ws.Cells.Item(89, 8).Value = 3004.5  # H89: 3019.45 -> 3004.5
$ws.Cells.Item(89, 9).Value = 2446.0667  # I89: 2368.125 -> 2446.0667
$ws.Cells.Item(89, 10).Value = 4679.8  # J89: 5624.75 -> 4679.8
$ws.Cells.Item(89, 11).Value = 12230.3335  # K89: 11840.625 -> 12230.3335
$ws.Cells.Item(89, 12).Value = 23399  # L89: 28123.75 -> 23399
$ws.Cells.Item(89, 13).Value = -6614.333499999999  # M89: -6224.625 -> -6614.333499999999
$ws.Cells.Item(89, 14).Value = -34631  # N89: -39355.75 -> -34631
$ws.Cells.Item(116, 8).Value = 15580.454  # H116: 13085.286 -> 15580.454
$ws.Cells.Item(116, 9).Value = 13948.125  # I116: 13259.4 -> 13948.125
$ws.Cells.Item(116, 10).Value = 19933.334  # J116: 12650 -> 19933.334
$ws.Cells.Item(116, 11).Value = 13948.125  # K116: 13259.4 -> 13948.125
$ws.Cells.Item(116, 12).Value = 19933.334  # L116: 12650 -> 19933.334
$ws.Cells.Item(116, 13).Value = -10506.125  # M116: -9817.4 -> -10506.125
$ws.Cells.Item(116, 14).Value = -26817.334  # N116: -19534 -> -26817.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2741.258  # H32: 2508.7646 -> 2741.258
$ws.Cells.Item(32, 9).Value = 2330.3276  # I32: 2121.8281 -> 2330.3276
$ws.Cells.Item(32, 11).Value = 2330.3276  # K32: 2121.8281 -> 2330.3276
$ws.Cells.Item(32, 13).Value = -2043.3276  # M32: -1834.8281 -> -2043.3276
$ws.Cells.Item(108, 8).Value = 45000  # H108: 49000 -> 45000
$ws.Cells.Item(108, 9).Value = 45000  # I108: 49000 -> 45000
$ws.Cells.Item(108, 11).Value = 45000  # K108: 49000 -> 45000
$ws.Cells.Item(108, 13).Value = -41160  # M108: -45160 -> -41160
$ws.Cells.Item(122, 8).Value = 2132.7778  # H122: 2167 -> 2132.7778
$ws.Cells.Item(122, 9).Value = 1970.7142  # I122: 1989.3334 -> 1970.7142
$ws.Cells.Item(122, 11).Value = 5912.142599999999  # K122: 5968.0002 -> 5912.142599999999
$ws.Cells.Item(122, 13).Value = -3462.142599999999  # M122: -3518.0002 -> -3462.142599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 0  # H35: 89000.336 -> 0
$ws.Cells.Item(35, 9).Value = 0  # I35: 67000 -> 0
$ws.Cells.Item(35, 10).Value = 0  # J35: 100000.5 -> 0
$ws.Cells.Item(35, 11).Value = 0  # K35: 67000 -> 0
$ws.Cells.Item(35, 12).ClearContents()  # L35: 100000.5 -> (removed)
$ws.Cells.Item(35, 13).ClearContents()  # M35: -66690 -> (removed)
$ws.Cells.Item(35, 14).Value = 0  # N35: -100620.5 -> 0
$ws.Cells.Item(75, 8).Value = 15214  # H75: 9214 -> 15214
$ws.Cells.Item(75, 9).Value = 15214  # I75: 9214 -> 15214
$ws.Cells.Item(75, 11).Value = 15214  # K75: 9214 -> 15214
$ws.Cells.Item(75, 13).Value = -14278  # M75: -8278 -> -14278
$ws.Cells.Item(78, 8).Value = 15214  # H78: 9214 -> 15214
$ws.Cells.Item(78, 9).Value = 15214  # I78: 9214 -> 15214
$ws.Cells.Item(78, 11).Value = 45642  # K78: 27642 -> 45642
$ws.Cells.Item(78, 13).Value = -40962  # M78: -22962 -> -40962
$ws.Cells.Item(105, 8).Value = 5746  # H105: 6116.5835 -> 5746
$ws.Cells.Item(105, 9).Value = 5105.4443  # I105: 5581.25 -> 5105.4443
$ws.Cells.Item(105, 11).Value = 5105.4443  # K105: 5581.25 -> 5105.4443
$ws.Cells.Item(105, 13).Value = -3358.4443  # M105: -3834.25 -> -3358.4443
$ws.Cells.Item(107, 8).Value = 1441.6842  # H107: 1524.5625 -> 1441.6842
$ws.Cells.Item(107, 9).Value = 1376.0588  # I107: 1426.2 -> 1376.0588
$ws.Cells.Item(107, 10).Value = 1999.5  # J107: 3000 -> 1999.5
$ws.Cells.Item(107, 11).Value = 1376.0588  # K107: 1426.2 -> 1376.0588
$ws.Cells.Item(107, 12).Value = 1999.5  # L107: 3000 -> 1999.5
$ws.Cells.Item(107, 13).Value = 543.9412  # M107: 493.8 -> 543.9412
$ws.Cells.Item(107, 14).Value = -5839.5  # N107: -6840 -> -5839.5
$ws.Cells.Item(129, 8).Value = 0  # H129: 67777 -> 0
$ws.Cells.Item(129, 10).Value = 0  # J129: 67777 -> 0
$ws.Cells.Item(129, 12).ClearContents()  # L129: 67777 -> (removed)
$ws.Cells.Item(129, 14).Value = 0  # N129: -77777 -> 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2683.5833  # H58: 2740.4 -> 2683.5833
$ws.Cells.Item(58, 10).Value = 2449.75  # J58: 2500 -> 2449.75
$ws.Cells.Item(58, 12).Value = 2449.75  # L58: 2500 -> 2449.75
$ws.Cells.Item(58, 14).Value = -2855.75  # N58: -2906 -> -2855.75
$ws.Cells.Item(82, 8).Value = 0  # H82: 19499.5 -> 0
$ws.Cells.Item(82, 10).Value = 0  # J82: 19499.5 -> 0
$ws.Cells.Item(82, 12).ClearContents()  # L82: 19499.5 -> (removed)
$ws.Cells.Item(82, 14).Value = 0  # N82: -20221.5 -> 0
$ws.Cells.Item(85, 8).Value = 0  # H85: 19499.5 -> 0
$ws.Cells.Item(85, 10).Value = 0  # J85: 19499.5 -> 0
$ws.Cells.Item(85, 12).ClearContents()  # L85: 19499.5 -> (removed)
$ws.Cells.Item(85, 14).Value = 0  # N85: -21995.5 -> 0
$ws.Cells.Item(136, 8).Value = 2683.5833  # H136: 2740.4 -> 2683.5833
$ws.Cells.Item(136, 10).Value = 2449.75  # J136: 2500 -> 2449.75
$ws.Cells.Item(136, 12).Value = 7349.25  # L136: 7500 -> 7349.25
$ws.Cells.Item(136, 14).Value = -12449.25  # N136: -12600 -> -12449.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(49, 8).Value = 0  # H49: 1000 -> 0
$ws.Cells.Item(49, 10).Value = 0  # J49: 1000 -> 0
$ws.Cells.Item(49, 12).ClearContents()  # L49: 3000 -> (removed)
$ws.Cells.Item(49, 14).Value = 0  # N49: -3312 -> 0
$ws.Cells.Item(63, 8).Value = 18199.334  # H63: 11319.6 -> 18199.334
$ws.Cells.Item(63, 9).Value = 5000  # I63: 2333.3333 -> 5000
$ws.Cells.Item(63, 11).Value = 15000  # K63: 6999.999899999999 -> 15000
$ws.Cells.Item(63, 13).Value = -14251  # M63: -6250.999899999999 -> -14251
$ws.Cells.Item(66, 8).Value = 18199.334  # H66: 11319.6 -> 18199.334
$ws.Cells.Item(66, 9).Value = 5000  # I66: 2333.3333 -> 5000
$ws.Cells.Item(66, 11).Value = 45000  # K66: 20999.9997 -> 45000
$ws.Cells.Item(66, 13).Value = -41256  # M66: -17255.9997 -> -41256
$ws.Cells.Item(86, 8).Value = 382.22223  # H86: 618.1429000000001 -> 382.22223
$ws.Cells.Item(86, 9).Value = 311.2  # I86: 361 -> 311.2
$ws.Cells.Item(86, 10).Value = 471  # J86: 961 -> 471
$ws.Cells.Item(86, 11).Value = 933.5999999999999  # K86: 1083 -> 933.5999999999999
$ws.Cells.Item(86, 12).Value = 1413  # L86: 2883 -> 1413
$ws.Cells.Item(86, 13).Value = 252.4000000000001  # M86: 103 -> 252.4000000000001
$ws.Cells.Item(86, 14).Value = -3785  # N86: -5255 -> -3785
$ws.Cells.Item(87, 8).Value = 6512.4614  # H87: 6962.5386 -> 6512.4614
$ws.Cells.Item(87, 9).Value = 4238.2  # I87: 4677 -> 4238.2
$ws.Cells.Item(87, 10).Value = 14093.333  # J87: 12105 -> 14093.333
$ws.Cells.Item(87, 11).Value = 12714.6  # K87: 14031 -> 12714.6
$ws.Cells.Item(87, 12).Value = 42279.999  # L87: 36315 -> 42279.999
$ws.Cells.Item(87, 13).Value = -11466.6  # M87: -12783 -> -11466.6
$ws.Cells.Item(87, 14).Value = -44775.999  # N87: -38811 -> -44775.999
$ws.Cells.Item(89, 8).Value = 382.22223  # H89: 618.1429000000001 -> 382.22223
$ws.Cells.Item(89, 9).Value = 311.2  # I89: 361 -> 311.2
$ws.Cells.Item(89, 10).Value = 471  # J89: 961 -> 471
$ws.Cells.Item(89, 11).Value = 2800.8  # K89: 3249 -> 2800.8
$ws.Cells.Item(89, 12).Value = 4239  # L89: 8649 -> 4239
$ws.Cells.Item(89, 13).Value = 3127.2  # M89: 2679 -> 3127.2
$ws.Cells.Item(89, 14).Value = -16095  # N89: -20505 -> -16095
$ws.Cells.Item(90, 8).Value = 6512.4614  # H90: 6962.5386 -> 6512.4614
$ws.Cells.Item(90, 9).Value = 4238.2  # I90: 4677 -> 4238.2
$ws.Cells.Item(90, 10).Value = 14093.333  # J90: 12105 -> 14093.333
$ws.Cells.Item(90, 11).Value = 38143.8  # K90: 42093 -> 38143.8
$ws.Cells.Item(90, 12).Value = 126839.997  # L90: 108945 -> 126839.997
$ws.Cells.Item(90, 13).Value = -31903.8  # M90: -35853 -> -31903.8
$ws.Cells.Item(90, 14).Value = -139319.997  # N90: -121425 -> -139319.997

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 950  # H7: 100000 -> 950
$ws.Cells.Item(7, 9).Value = 900  # I7: 0 -> 900
$ws.Cells.Item(7, 10).Value = 1000  # J7: 100000 -> 1000
$ws.Cells.Item(7, 11).Value = 900  # K7: 0 -> 900
$ws.Cells.Item(7, 12).Value = 1000  # L7: 100000 -> 1000
$ws.Cells.Item(7, 13).Value = -788  # M7: None -> -788
$ws.Cells.Item(7, 14).Value = -1224  # N7: -100224 -> -1224
$ws.Cells.Item(8, 8).Value = 950  # H8: 100000 -> 950
$ws.Cells.Item(8, 9).Value = 900  # I8: 0 -> 900
$ws.Cells.Item(8, 10).Value = 1000  # J8: 100000 -> 1000
$ws.Cells.Item(8, 11).Value = 900  # K8: 0 -> 900
$ws.Cells.Item(8, 12).Value = 1000  # L8: 100000 -> 1000
$ws.Cells.Item(8, 13).Value = -761  # M8: None -> -761
$ws.Cells.Item(8, 14).Value = -1278  # N8: -100278 -> -1278
$ws.Cells.Item(11, 8).Value = 5714714  # H11: 7525250 -> 5714714
$ws.Cells.Item(11, 9).Value = 6667000  # I11: 10033333 -> 6667000
$ws.Cells.Item(11, 10).Value = 999  # J11: 1000 -> 999
$ws.Cells.Item(11, 11).Value = 6667000  # K11: 10033333 -> 6667000
$ws.Cells.Item(11, 12).Value = 999  # L11: 1000 -> 999
$ws.Cells.Item(11, 13).Value = -6666861  # M11: -10033194 -> -6666861
$ws.Cells.Item(11, 14).Value = -1277  # N11: -1278 -> -1277
$ws.Cells.Item(54, 8).Value = 18123.75  # H54: 12148.667 -> 18123.75
$ws.Cells.Item(54, 9).Value = 7500  # I54: 0 -> 7500
$ws.Cells.Item(54, 10).Value = 21665  # J54: 12148.667 -> 21665
$ws.Cells.Item(54, 11).Value = 7500  # K54: 0 -> 7500
$ws.Cells.Item(54, 12).Value = 21665  # L54: 12148.667 -> 21665
$ws.Cells.Item(54, 13).Value = -7110  # M54: None -> -7110
$ws.Cells.Item(54, 14).Value = -22445  # N54: -12928.667 -> -22445
$ws.Cells.Item(132, 8).Value = 225537.64  # H132: 235932.42 -> 225537.64
$ws.Cells.Item(132, 9).Value = 225537.64  # I132: 235932.42 -> 225537.64
$ws.Cells.Item(132, 11).Value = 676612.92  # K132: 707797.26 -> 676612.92
$ws.Cells.Item(132, 13).Value = -674082.92  # M132: -705267.26 -> -674082.92

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 2261.2  # H46: 2393.8333 -> 2261.2
$ws.Cells.Item(46, 9).Value = 2337.25  # I46: 2374.75 -> 2337.25
$ws.Cells.Item(46, 10).Value = 2174.2856  # J46: 2432 -> 2174.2856
$ws.Cells.Item(46, 11).Value = 2337.25  # K46: 2374.75 -> 2337.25
$ws.Cells.Item(46, 12).Value = 2174.2856  # L46: 2432 -> 2174.2856
$ws.Cells.Item(46, 13).Value = -2149.25  # M46: -2186.75 -> -2149.25
$ws.Cells.Item(46, 14).Value = -2550.2856  # N46: -2808 -> -2550.2856
$ws.Cells.Item(61, 8).Value = 4082.4167  # H61: 4261.619 -> 4082.4167
$ws.Cells.Item(61, 9).Value = 3999.111  # I61: 4233.2 -> 3999.111
$ws.Cells.Item(61, 10).Value = 4332.3335  # J61: 4332.6665 -> 4332.3335
$ws.Cells.Item(61, 11).Value = 3999.111  # K61: 4233.2 -> 3999.111
$ws.Cells.Item(61, 12).Value = 4332.3335  # L61: 4332.6665 -> 4332.3335
$ws.Cells.Item(61, 13).Value = -3797.111  # M61: -4031.2 -> -3797.111
$ws.Cells.Item(61, 14).Value = -4736.3335  # N61: -4736.6665 -> -4736.3335
$ws.Cells.Item(68, 8).Value = 0  # H68: 1219.4 -> 0
$ws.Cells.Item(68, 9).Value = 0  # I68: 1173.5 -> 0
$ws.Cells.Item(68, 10).Value = 0  # J68: 1250 -> 0
$ws.Cells.Item(68, 11).Value = 0  # K68: 1173.5 -> 0
$ws.Cells.Item(68, 12).ClearContents()  # L68: 1250 -> (removed)
$ws.Cells.Item(68, 13).ClearContents()  # M68: -424.5 -> (removed)
$ws.Cells.Item(68, 14).Value = 0  # N68: -2748 -> 0
$ws.Cells.Item(71, 8).Value = 0  # H71: 1219.4 -> 0
$ws.Cells.Item(71, 9).Value = 0  # I71: 1173.5 -> 0
$ws.Cells.Item(71, 10).Value = 0  # J71: 1250 -> 0
$ws.Cells.Item(71, 11).Value = 0  # K71: 5867.5 -> 0
$ws.Cells.Item(71, 12).ClearContents()  # L71: 6250 -> (removed)
$ws.Cells.Item(71, 13).ClearContents()  # M71: -2123.5 -> (removed)
$ws.Cells.Item(71, 14).Value = 0  # N71: -13738 -> 0
$ws.Cells.Item(82, 8).Value = 2486.7646  # H82: 2523.4375 -> 2486.7646
$ws.Cells.Item(82, 9).Value = 2227.2222  # I82: 2268.125 -> 2227.2222
$ws.Cells.Item(82, 11).Value = 2227.2222  # K82: 2268.125 -> 2227.2222
$ws.Cells.Item(82, 13).Value = -1866.2222  # M82: -1907.125 -> -1866.2222
$ws.Cells.Item(85, 8).Value = 2486.7646  # H85: 2523.4375 -> 2486.7646
$ws.Cells.Item(85, 9).Value = 2227.2222  # I85: 2268.125 -> 2227.2222
$ws.Cells.Item(85, 11).Value = 2227.2222  # K85: 2268.125 -> 2227.2222
$ws.Cells.Item(85, 13).Value = -979.2222000000002  # M85: -1020.125 -> -979.2222000000002
$ws.Cells.Item(100, 8).Value = 62102.055  # H100: 62263.223 -> 62102.055
$ws.Cells.Item(100, 9).Value = 67427.06  # I100: 71788.92999999999 -> 67427.06
$ws.Cells.Item(100, 10).Value = 19502  # J100: 14634.667 -> 19502
$ws.Cells.Item(100, 11).Value = 67427.06  # K100: 71788.92999999999 -> 67427.06
$ws.Cells.Item(100, 12).Value = 19502  # L100: 14634.667 -> 19502
$ws.Cells.Item(100, 13).Value = -66886.06  # M100: -71247.92999999999 -> -66886.06
$ws.Cells.Item(100, 14).Value = -20584  # N100: -15716.667 -> -20584
$ws.Cells.Item(113, 8).Value = 4082.4167  # H113: 4261.619 -> 4082.4167
$ws.Cells.Item(113, 9).Value = 3999.111  # I113: 4233.2 -> 3999.111
$ws.Cells.Item(113, 10).Value = 4332.3335  # J113: 4332.6665 -> 4332.3335
$ws.Cells.Item(113, 11).Value = 3999.111  # K113: 4233.2 -> 3999.111
$ws.Cells.Item(113, 12).Value = 4332.3335  # L113: 4332.6665 -> 4332.3335
$ws.Cells.Item(113, 13).Value = -1829.111  # M113: -2063.2 -> -1829.111
$ws.Cells.Item(113, 14).Value = -8672.333500000001  # N113: -8672.666499999999 -> -8672.333500000001
$ws.Cells.Item(132, 8).Value = 6135.8184  # H132: 6049.3 -> 6135.8184
$ws.Cells.Item(132, 9).Value = 4999.1665  # I132: 5332.3335 -> 4999.1665
$ws.Cells.Item(132, 10).Value = 7499.8  # J132: 7124.75 -> 7499.8
$ws.Cells.Item(132, 11).Value = 14997.4995  # K132: 15997.0005 -> 14997.4995
$ws.Cells.Item(132, 12).Value = 22499.4  # L132: 21374.25 -> 22499.4
$ws.Cells.Item(132, 13).Value = -12467.4995  # M132: -13467.0005 -> -12467.4995
$ws.Cells.Item(132, 14).Value = -27559.4  # N132: -26434.25 -> -27559.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(38, 8).Value = 8750  # H38: 9327.5 -> 8750
$ws.Cells.Item(38, 9).Value = 0  # I38: 14750 -> 0
$ws.Cells.Item(38, 10).Value = 8750  # J38: 7520 -> 8750
$ws.Cells.Item(38, 11).Value = 0  # K38: 14750 -> 0
$ws.Cells.Item(38, 12).ClearContents()  # L38: 7520 -> (removed)
$ws.Cells.Item(38, 13).Value = 8750  # M38: -14277 -> 8750
$ws.Cells.Item(38, 14).Value = -9696  # N38: -8466 -> -9696
$ws.Cells.Item(51, 8).Value = 0  # H51: 39999 -> 0
$ws.Cells.Item(51, 9).Value = 0  # I51: 39999 -> 0
$ws.Cells.Item(51, 11).Value = 0  # K51: 39999 -> 0
$ws.Cells.Item(51, 13).ClearContents()  # M51: -39489 -> (removed)
$ws.Cells.Item(86, 8).Value = 50324.5  # H86: 50325 -> 50324.5
$ws.Cells.Item(86, 10).Value = 50324.5  # J86: 50325 -> 50324.5
$ws.Cells.Item(86, 12).Value = 50324.5  # L86: 50325 -> 50324.5
$ws.Cells.Item(86, 14).Value = -52570.5  # N86: -52571 -> -52570.5
$ws.Cells.Item(89, 8).Value = 50324.5  # H89: 50325 -> 50324.5
$ws.Cells.Item(89, 10).Value = 50324.5  # J89: 50325 -> 50324.5
$ws.Cells.Item(89, 12).Value = 251622.5  # L89: 251625 -> 251622.5
$ws.Cells.Item(89, 14).Value = -262854.5  # N89: -262857 -> -262854.5
$ws.Cells.Item(100, 8).Value = 965.3684  # H100: 856.0968 -> 965.3684
$ws.Cells.Item(100, 9).Value = 791.3103599999999  # I100: 838.8148 -> 791.3103599999999
$ws.Cells.Item(100, 10).Value = 1526.2222  # J100: 972.75 -> 1526.2222
$ws.Cells.Item(100, 11).Value = 1582.62072  # K100: 1677.6296 -> 1582.62072
$ws.Cells.Item(100, 12).Value = 3052.4444  # L100: 1945.5 -> 3052.4444
$ws.Cells.Item(100, 13).Value = -1041.62072  # M100: -1136.6296 -> -1041.62072
$ws.Cells.Item(100, 14).Value = -4134.4444  # N100: -3027.5 -> -4134.4444
$ws.Cells.Item(109, 8).Value = 39933.332  # H109: 39950 -> 39933.332
$ws.Cells.Item(109, 10).Value = 39933.332  # J109: 39950 -> 39933.332
$ws.Cells.Item(109, 12).Value = 39933.332  # L109: 39950 -> 39933.332
$ws.Cells.Item(109, 14).Value = -42707.332  # N109: -42724 -> -42707.332
$ws.Cells.Item(136, 8).Value = 2273.318  # H136: 2226.5652 -> 2273.318
$ws.Cells.Item(136, 9).Value = 2001  # I136: 1960.85 -> 2001
$ws.Cells.Item(136, 11).Value = 6003  # K136: 5882.549999999999 -> 6003
$ws.Cells.Item(136, 13).Value = -3453  # M136: -3332.549999999999 -> -3453
